# tracklist.xlsx edit: mkiwanuka2 (Michael Kiwanuka tracklist) -> maxwell4
# (Maxwell - blackSUMMERS'night tracklist), adds a 12th track.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------------
# 1) Track data (Sheet1 and Sheet3 are identical copies of the web-query
#    result range, rows 2..13 = tracks 1..12).
# ---------------------------------------------------------------------------

$titles = @(
    "All the Ways Love Can Feel",
    "The Fall",
    "III",
    "Lake by the Ocean",
    "Fingers Crossed",
    "Hostage",
    "1990X",
    "Gods",
    "Lost",
    "Of All Kind",
    "Listen Hear",
    "Night"
)

$composers = @(
    "Hod David / Musze / Travis Sayles",
    "Hod David / Shedrick Mitchell / Musze / Travis Sayles",
    "Hod David / Musze",
    "Hod David / Musze",
    "Hod David / Musze",
    "Hod David / Musze",
    "Hod David / Musze",
    "Hod David / Musze",
    "Stuart Matthewman / Musze",
    "Hod David / Musze",
    "Stuart Matthewman / Musze",
    "Earth"
)

$performer = "Maxwell"

$times = @(
    0.22291666666666665,
    0.17430555555555557,
    0.19930555555555554,
    0.16527777777777777,
    0.19305555555555554,
    0.16319444444444445,
    0.19722222222222222,
    0.14930555555555555,
    0.16458333333333333,
    0.15416666666666667,
    0.15138888888888888,
    0.015972222222222224
)

for ($sheetIdx = 0; $sheetIdx -lt 2; $sheetIdx++) {
    if ($sheetIdx -eq 0) { $ws = $ws1 } else { $ws = $ws3 }

    for ($i = 0; $i -lt 12; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $i + 1
        $ws.Cells.Item($row, 2).Value = $titles[$i]
        $ws.Cells.Item($row, 3).Value = $composers[$i]
        $ws.Cells.Item($row, 4).Value = $performer
        $ws.Cells.Item($row, 5).Value = $times[$i]
    }

    # Column widths (B..E) changed to fit the new, slightly different text.
    $ws.Columns.Item(2).ColumnWidth = 25.5546875
    $ws.Columns.Item(3).ColumnWidth = 50
    $ws.Columns.Item(4).ColumnWidth = 10
    $ws.Columns.Item(5).ColumnWidth = 9.44140625
}

# ---------------------------------------------------------------------------
# 2) Defined names: mkiwanuka2 -> maxwell4, range grows from row 12 to row 13
#    (set RefersTo before Name to dodge a rename/re-index quirk when both
#    sheet-scoped names end up sharing the same short name).
# ---------------------------------------------------------------------------

$n1 = $wb.Names.Item(1)
$n2 = $wb.Names.Item(2)
$n1.RefersTo = "=Sheet1!`$A`$1:`$E`$13"
$n2.RefersTo = "=Sheet3!`$A`$1:`$E`$13"
$n1.Name = "maxwell4"
$n2.Name = "maxwell4"

# ---------------------------------------------------------------------------
# 3) Sheet2 (the formatted "pretty print" sheet) is 100% formulas driven off
#    Sheet1, so it recalculates on its own. Only the view's selection needs
#    nudging to cover the two extra rows now in use.
# ---------------------------------------------------------------------------

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("A3:K16").Select()
